# Append: 2025-09-28 18:28 JST
# Update the "取得日時" (acquired timestamp) column for the existing data
# rows from 2025-09-28 18:22:42 to 2025-09-28 18:28:38.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-09-28 18:28:38"

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
